# Applies the edits described by the diff:
#  - A1 header "Gen" -> "MaxFES"
#  - Column A values (generation counts) -> MaxFES fractions
#  - Remove the extraneous "Run 50" column (AZ) entirely, shifting the old
#    "Mean" column (BA) left into AZ
#  - Recompute the Mean column (now AZ) as the average of the 50 run
#    columns (B:AY), since the previous mean included the removed run

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Update header and A column values
$ws.Range("A1").Value = "MaxFES"

$fesValues = @(0, 0.001, 0.01, 0.1, 0.2, 0.3, 0.4, 0.5, 0.6, 0.7, 0.8, 0.9, 1)
for ($i = 0; $i -lt $fesValues.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $fesValues[$i]
}

# 2. Remove the "Run 50" column (AZ). This shifts the old "Mean" column
#    (BA) left to become the new AZ column, and the sheet dimension
#    shrinks from BA14 to AZ14 automatically.
$ws.Columns("AZ:AZ").Delete()

# 3. Recompute the Mean column (now AZ) header + values since the
#    removed run is no longer part of the average.
$ws.Range("AZ1").Value = "Mean"

$meanValues = @(
    134.84050197,
    123.79930704,
    80.64210699,
    26.99823101,
    18.12991606,
    13.96441838,
    11.08216914,
    9.98184403,
    8.62776201,
    7.48715838,
    6.73186696,
    6.23368929,
    5.74854146
)
for ($i = 0; $i -lt $meanValues.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 52).Value = $meanValues[$i]
}
